# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 33-34) above the existing Higo
# price history, pushing the previous rows 33-44 down to rows 35-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 33:34 - everything currently at/after row 33
# shifts down by two rows (matches the new dimension A1:T46).
$ws.Range("A33:T34").EntireRow.Insert()

# Populate the newly inserted row 33 with this week's "Primera" record.
$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 45027
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100101
$ws.Range("H33").Value = "Berries"
$ws.Range("I33").Value = 100101006
$ws.Range("J33").Value = "Higo"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 20000
$ws.Range("Q33").Value = "$/bandeja 7 kilos"
$ws.Range("R33").Value = "Región Metropolitana"
$ws.Range("S33").Value = 2857
$ws.Range("T33").Value = 7

# Populate the newly inserted row 34 with this week's "Segunda" record.
$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 45027
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100101
$ws.Range("H34").Value = "Berries"
$ws.Range("I34").Value = 100101006
$ws.Range("J34").Value = "Higo"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Segunda"
$ws.Range("M34").Value = 80
$ws.Range("N34").Value = 14000
$ws.Range("O34").Value = 14000
$ws.Range("P34").Value = 14000
$ws.Range("Q34").Value = "$/bandeja 7 kilos"
$ws.Range("R34").Value = "Región Metropolitana"
$ws.Range("S34").Value = 2000
$ws.Range("T34").Value = 7
